$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: ManageQAData -- path (C) entered before name (B)
$ws.Cells.Item(13, 1).Value = "test"
$ws.Cells.Item(13, 3).Value = "\Testdata\Non_Oncology\DataFiles\ManageQAData\ManageQADataPage_NonOnco_Data.xlsx"
$ws.Cells.Item(13, 2).Value = "nononcology_manageqadata"

# Row 14: ManageExcludedPublications -- path (C) entered before name (B)
$ws.Cells.Item(14, 1).Value = "test"
$ws.Cells.Item(14, 3).Value = "\Testdata\Non_Oncology\DataFiles\ManageExcludedPublications\ExcludedPubPage_Data.xlsx"
$ws.Cells.Item(14, 2).Value = "nononcology_manageexcludedpub"

# Row 15: PRISMA -- name (B) entered before path (C)
$ws.Cells.Item(15, 1).Value = "test"
$ws.Cells.Item(15, 2).Value = "nononcology_prisma"
$ws.Cells.Item(15, 3).Value = "\Testdata\Non_Oncology\DataFiles\Protocol_Page\PRISMA\PRISMA_Nononco_Data.xlsx"

# Row 16: PICOS -- name (B) entered before path (C)
$ws.Cells.Item(16, 1).Value = "test"
$ws.Cells.Item(16, 2).Value = "nononcology_picos"
$ws.Cells.Item(16, 3).Value = "\Testdata\Non_Oncology\DataFiles\Protocol_Page\PICOS\PICOS_Nononco_Data.xlsx"

# Row 17: SearchStrategy -- name (B) entered before path (C)
$ws.Cells.Item(17, 1).Value = "test"
$ws.Cells.Item(17, 2).Value = "nononcology_searchstrategy"
$ws.Cells.Item(17, 3).Value = "\Testdata\Non_Oncology\DataFiles\Protocol_Page\SearchStrategy\SearchStrategy_Nononco_Data.xlsx"

$ws.Range("B17").Select() | Out-Null

$wb.Save()
